# feat: fix 2 issue
#  1) Insert a new "Kasbon" column header before the existing
#     "Total Penghasilan" column - this pushes "Total Penghasilan",
#     "Total Pengurangan" and "Penerimaan Bersih" one column to the
#     right (I->J, J->K, K->L) and adds the new header in column I.
#  2) The saved selection moves to I2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the header cells one column to the right, starting from the
# rightmost one so we never overwrite a source cell before it has been
# copied. Copy() carries over both the value and the cell formatting
# (bold + centered header style), which matches what Excel does when a
# column is inserted.
$ws.Range("K1").Copy($ws.Range("L1"))
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("I1").Copy($ws.Range("J1"))

# New column header, reusing the (now-duplicated) header style already
# sitting in I1.
$ws.Range("I1").Value2 = "Kasbon"

# Match the saved selection recorded in the workbook.
[void]$ws.Range("I2").Select()
